# ------------------------------------------------------------------
# Language.xlsx restructuring:
#   - rename "Sheet1" -> "Comm" and extend/relabel its data rows
#   - add four new sheets: Property, Guild, Tip, Item
#   - populate each with the localisation strings from the new layout
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ===================== Comm (was Sheet1) =====================
$comm = $wb.Worksheets.Item(1)
$comm.Name = "Comm"

# Resize columns to the new, wider layout
$comm.Columns.Item(1).ColumnWidth = 31.16
$comm.Columns.Item(2).ColumnWidth = 23.79
$comm.Columns.Item(3).ColumnWidth = 22.29

# Rows 2-7: column A gets new "Langage_Comm_n" ids, column C gets the
# new Chinese strings; column B keeps its existing Langage_n values.
$comm.Range("A2").Value = "Langage_Comm_1"
$comm.Range("C2").Value = "确认"

$comm.Range("A3").Value = "Langage_Comm_2"
$comm.Range("C3").Value = "取消"

$comm.Range("A4").Value = "Langage_Comm_3"
$comm.Range("C4").Value = "登录"

$comm.Range("A5").Value = "Langage_Comm_4"
$comm.Range("C5").Value = "创建角色"

$comm.Range("A6").Value = "Langage_Comm_5"
$comm.Range("C6").Value = "进入游戏"

$comm.Range("A7").Value = "Langage_Comm_6"
$comm.Range("C7").Value = "中文_6"

# Rows 8-12: blank filler rows, formatted like row 2
$comm.Range("A2:C2").Copy()
$comm.Range("A8:C12").PasteSpecial($xlPasteFormats)
$comm.Range("A8:C12").ClearContents()

$comm.Range("C8").Select()

# ===================== Property =====================
$property = $wb.Worksheets.Add([Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$property.Name = "Property"

$property.Columns.Item(1).ColumnWidth = 50.54

$comm.Range("A1:C1").Copy()
$property.Range("A1:C1").PasteSpecial($xlPasteFormats)

$comm.Range("B1").Copy()
$property.Range("A1:A28").PasteSpecial($xlPasteFormats)

$property.Range("A1").Value = "ID"
$property.Range("B1").Value = "English"
$property.Range("C1").Value = "Chinese"

$property.Range("A2").Value = "Langage_HP"
$property.Range("A3").Value = "Langage_MAXHP"
$property.Range("A4").Value = "Langage_MP"
$property.Range("A5").Value = "Langage_MAXMP"
$property.Range("A6").Value = "Langage_VP"
$property.Range("A7").Value = "Langage_ATTACK"

$property.Rows.Item(1).Select()

# ===================== Guild =====================
$guild = $wb.Worksheets.Add([Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$guild.Name = "Guild"

$comm.Range("A1:C1").Copy()
$guild.Range("A1:C1").PasteSpecial($xlPasteFormats)

$guild.Columns.Item(1).ColumnWidth = 31.16
$guild.Columns.Item(2).ColumnWidth = 23.79
$guild.Columns.Item(3).ColumnWidth = 22.29

$guild.Range("A1").Value = "ID"
$guild.Range("B1").Value = "English"
$guild.Range("C1").Value = "Chinese"

$comm.Range("A2:C2").Copy()
$guild.Range("A2:C12").PasteSpecial($xlPasteFormats)
$guild.Range("A3:C12").ClearContents()

$comm.Range("A2:C2").Copy()
$guild.Range("A16:C16").PasteSpecial($xlPasteFormats)
$guild.Range("A16:C16").ClearContents()

$guild.Range("A2").Value = "Langage_Guild_1"
$guild.Range("B2").Value = "Langage_1"
$guild.Range("C2").Value = "确认要加入这个公会吗？点击确认加入"

$comm.Range("B1").Copy()
$guild.Range("A22").PasteSpecial($xlPasteFormats)
$guild.Range("A22").ClearContents()

$guild.Range("A12").Select()

# ===================== Tip =====================
$tip = $wb.Worksheets.Add([Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$tip.Name = "Tip"

$comm.Range("A1:C1").Copy()
$tip.Range("A1:C1").PasteSpecial($xlPasteFormats)

$tip.Range("A1").Value = "ID"
$tip.Range("B1").Value = "English"
$tip.Range("C1").Value = "Chinese"

$tip.Rows.Item(1).Select()

# ===================== Item =====================
$item = $wb.Worksheets.Add([Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$item.Name = "Item"

$comm.Range("A1:C1").Copy()
$item.Range("A1:C1").PasteSpecial($xlPasteFormats)

$item.Range("A1").Value = "ID"
$item.Range("B1").Value = "English"
$item.Range("C1").Value = "Chinese"

$item.Rows.Item(1).Select()

# Leave the workbook with Comm as the active/visible sheet.
$comm.Activate()
